# Refresh the cryptos price/volume snapshot (and fix two swapped-row pairs)
# to match the latest GitHub Actions scrape.
# Note: Price values that look like plain decimals (e.g. "1.00", "7.51")
# are written with a leading apostrophe so Excel keeps them as text
# instead of silently converting them to numbers (which would drop
# trailing zeros such as in "1.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '71.300.91'
$ws.Cells.Item(2, 5).Value = '  +0.52%  '
$ws.Cells.Item(3, 4).Value = '3.808.34'
$ws.Cells.Item(3, 5).Value = '  -1.04%  '
$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '''702.27'
$ws.Cells.Item(5, 5).Value = '  -0.68%  '
$ws.Cells.Item(6, 4).Value = '''171.50'
$ws.Cells.Item(6, 5).Value = '  -0.56%  '
$ws.Cells.Item(7, 4).Value = '3.808.09'
$ws.Cells.Item(7, 5).Value = '  -0.96%  '
$ws.Cells.Item(8, 5).Value = '  +0.08%  '
$ws.Cells.Item(9, 5).Value = '  +0.45%  '
$ws.Cells.Item(10, 5).Value = '  -0.98%  '
$ws.Cells.Item(11, 4).Value = '''7.51'
$ws.Cells.Item(11, 5).Value = '  +2.71%  '
$ws.Cells.Item(12, 5).Value = '  +3.74%  '
$ws.Cells.Item(13, 4).Value = '''0.0000252'
$ws.Cells.Item(13, 5).Value = '  -1.76%  '
$ws.Cells.Item(14, 4).Value = '''35.99'
$ws.Cells.Item(14, 5).Value = '  -1.92%  '
$ws.Cells.Item(15, 4).Value = '4.448.31'
$ws.Cells.Item(15, 5).Value = '  -1.06%  '
$ws.Cells.Item(16, 4).Value = '3.826.19'
$ws.Cells.Item(16, 5).Value = '  -0.97%  '
$ws.Cells.Item(17, 4).Value = '71.305.28'
$ws.Cells.Item(17, 5).Value = '  +0.48%  '
$ws.Cells.Item(18, 4).Value = '''7.18'
$ws.Cells.Item(18, 5).Value = '  -0.14%  '
$ws.Cells.Item(19, 5).Value = '  -0.50%  '
$ws.Cells.Item(20, 4).Value = '''17.45'
$ws.Cells.Item(20, 5).Value = '  +0.55%  '
$ws.Cells.Item(21, 4).Value = '''513.79'
$ws.Cells.Item(21, 5).Value = '  +4.19%  '
$ws.Cells.Item(22, 5).Value = '  -1.50%  '
$ws.Cells.Item(23, 4).Value = '''0.715'
$ws.Cells.Item(23, 5).Value = '  -0.23%  '
$ws.Cells.Item(24, 4).Value = '''84.05'
$ws.Cells.Item(24, 5).Value = '  -1.46%  '
$ws.Cells.Item(25, 5).Value = '  -3.07%  '
$ws.Cells.Item(26, 4).Value = '''12.22'
$ws.Cells.Item(26, 5).Value = '  +0.57%  '
$ws.Cells.Item(27, 4).Value = '3.956.12'
$ws.Cells.Item(27, 5).Value = '  -1.06%  '
$ws.Cells.Item(28, 4).Value = '''10.34'
$ws.Cells.Item(28, 5).Value = '  -2.80%  '
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$ws.Cells.Item(30, 5).Value = '  -3.97%  '
$ws.Cells.Item(31, 4).Value = '''3.02'
$ws.Cells.Item(31, 5).Value = '  -5.35%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).Value = '''2.23'
$ws.Cells.Item(32, 5).Value = '  -1.50%  '
$ws.Cells.Item(33, 2).Value = 'NEARProtocol'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(33, 4).Value = '''7.33'
$ws.Cells.Item(33, 5).Value = '  -2.14%  '
$ws.Cells.Item(34, 4).Value = '''29.06'
$ws.Cells.Item(34, 5).Value = '  -1.36%  '
$ws.Cells.Item(35, 4).Value = '''0.173'
$ws.Cells.Item(35, 5).Value = '  -3.69%  '
$ws.Cells.Item(36, 4).Value = '''9.14'
$ws.Cells.Item(36, 5).Value = '  -0.08%  '
$ws.Cells.Item(37, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(37, 4).Value = '3.770.02'
$ws.Cells.Item(37, 5).Value = '  -0.88%  '
$ws.Cells.Item(38, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(38, 4).Value = '''0.999'
$ws.Cells.Item(38, 5).Value = '  -0.09%  '
$ws.Cells.Item(39, 2).Value = 'Filecoin'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(39, 4).Value = '''6.46'
$ws.Cells.Item(39, 5).Value = '  +6.90%  '
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).Value = '''0.101'
$ws.Cells.Item(40, 5).Value = '  -2.09%  '
$ws.Cells.Item(41, 4).Value = '''2.39'
$ws.Cells.Item(41, 5).Value = '  +0.88%  '
$ws.Cells.Item(42, 5).Value = '  -2.02%  '
$ws.Cells.Item(43, 4).Value = '''3.27'
$ws.Cells.Item(43, 5).Value = '  -1.71%  '
$ws.Cells.Item(44, 5).Value = '  -0.02%  '
$ws.Cells.Item(45, 4).Value = '''173.44'
$ws.Cells.Item(45, 5).Value = '  +6.50%  '
$ws.Cells.Item(46, 5).Value = '  +0.05%  '
$ws.Cells.Item(47, 5).Value = '  -2.30%  '
$ws.Cells.Item(48, 4).Value = '''49.49'
$ws.Cells.Item(48, 5).Value = '  +1.56%  '
$ws.Cells.Item(49, 4).Value = '''424.47'
$ws.Cells.Item(49, 5).Value = '  +1.91%  '
$ws.Cells.Item(50, 5).Value = '  -0.43%  '
$ws.Cells.Item(51, 5).Value = '  -1.12%  '
